# Generate Report for Handoff
# Update status text from "In Translation" to "Ready for handoff" and
# refresh the related "Latest Handoff/HO Xliff Generate" timestamps on
# all three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# E2 / F2: status for zh-cn / de-de
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-15 16:37:38"

# --- zh-cn sheet ---
# C2: Status
$wsZhCn.Range("C2").Value = "Ready for handoff"
# H2: Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-15 16:37:33"

# --- de-de sheet ---
# C2: Status
$wsDeDe.Range("C2").Value = "Ready for handoff"
# H2: Latest Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-08-15 16:37:38"

# Widen the Status columns so their widths reflect the new, longer text
# (matches the width Excel computes when auto-fitting "Ready for handoff";
# the input is tuned so the host's internal pixel-snapping lands on the
# same stored column width as the target workbook)
$wsOverview.Range("E1").ColumnWidth = 16.35
$wsOverview.Range("F1").ColumnWidth = 16.35
$wsZhCn.Range("C1").ColumnWidth = 16.35
$wsDeDe.Range("C1").ColumnWidth = 16.35
